# Applies the "Added Done Point enumeration and Test Sequence Template" edit:
#  - Duplicates the header row (row 7) of the "Component Comparision" sheet
#    down to row 16 (a fresh copy of the Test Sequence Template header).
#  - Moves the active tab / selection from "Major Components " (cell B7:E7)
#    over to "Component Comparision" (cell D25), with "Major Components "
#    leaving its own cursor parked at J11.

$wb = $excel.ActiveWorkbook

$wsMajor = $wb.Worksheets.Item("Major Components ")
$wsComp  = $wb.Worksheets.Item("Component Comparision")

# --- Component Comparision: copy the header row (row 7) down to row 16 ---
$wsComp.Range("A7:O7").Copy($wsComp.Range("A16:O16"))
$wsComp.Rows.Item(16).RowHeight = $wsComp.Rows.Item(7).RowHeight

# Match the portrait page setup already used on "Major Components "
$wsComp.PageSetup.Orientation = 1

# --- Selections on each sheet ---
$wsMajor.Range("J11").Select()
$wsComp.Range("D25").Select()

# --- Make "Component Comparision" the active/visible tab ---
$wsComp.Activate()
